$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item(1)
$zhcn     = $wb.Worksheets.Item(2)
$dede     = $wb.Worksheets.Item(3)

$handedBackStatus = "Handed back: in sync with en-US"

# --- Update status text (shared across Overview E2/F2 and both table C2 cells) ---
$overview.Range("E2").Value = $handedBackStatus
$overview.Range("F2").Value = $handedBackStatus
$zhcn.Range("C2").Value = $handedBackStatus
$dede.Range("C2").Value = $handedBackStatus

# --- zh-cn sheet: fill in Latest Target File / Latest Handback File / Latest Handback DateTime ---
$zhcnI2 = $zhcn.Range("I2")
$zhcnI2.Hyperlinks.Add($zhcnI2, "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/95e08786cb0e0088b0a709705a7a425106ddbde1/e2e/7d71bcde-2188-4dea-9990-360038628121.md", "", "", "7d71bcde-2188-4dea-9990-360038628121.md")
$zhcn.Range("J2").Value = "7d71bcde-2188-4dea-9990-360038628121.d4ae73b8332663ba5dcae6116e1de887174e2bbf.zh-cn.xlf"
$zhcn.Range("K2").Value = "2016-09-02 07:08:00"

# --- de-de sheet: fill in Latest Target File / Latest Handback File / Latest Handback DateTime ---
$dedeI2 = $dede.Range("I2")
$dedeI2.Hyperlinks.Add($dedeI2, "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/95e08786cb0e0088b0a709705a7a425106ddbde1/e2e/7d71bcde-2188-4dea-9990-360038628121.md", "", "", "7d71bcde-2188-4dea-9990-360038628121.md")
$dede.Range("J2").Value = "7d71bcde-2188-4dea-9990-360038628121.d4ae73b8332663ba5dcae6116e1de887174e2bbf.de-de.xlf"
$dede.Range("K2").Value = "2016-09-02 07:08:15"

# --- widen columns to fit the new/longer content ---
# (ColumnWidth is specified in characters; the stored OOXML width = ColumnWidth + 5/6)
$overview.Range("E1").ColumnWidth = 29.166666666666668
$overview.Range("F1").ColumnWidth = 29.166666666666668

$zhcn.Range("C1").ColumnWidth = 29.166666666666668
$zhcn.Range("I1").ColumnWidth = 39.166666666666664
$zhcn.Range("J1").ColumnWidth = 39.166666666666664

$dede.Range("C1").ColumnWidth = 29.166666666666668
$dede.Range("I1").ColumnWidth = 39.166666666666664
$dede.Range("J1").ColumnWidth = 39.166666666666664

Write-Host "done"
